# "cambios carlos documentos y texto"
# The only substantive content edit in this commit: cell B6 on "Hoja1"
# (shared string "*") is changed to "**".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "**"

# The author also left the selection/scroll position at B6 (previously A7)
# when they saved the file.
$ws.Range("B6").Select()
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
